# Updates cryptos list figures (Price + Volume(1h)) for rows 2-51 of Sheet1,
# matching the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.217.82"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "2.267.76"
$ws.Range("E3").Value = "  -0.90%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.56"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.17"
$ws.Range("E6").Value = "  +1.48%  "
$ws.Range("E7").Value = "  -1.16%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.491"
$ws.Range("E9").Value = "  -0.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.15"
$ws.Range("E10").Value = "  -2.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0790"
$ws.Range("E11").Value = "  -1.69%  "
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.86"
$ws.Range("E13").Value = "  +1.59%  "
$ws.Range("D14").Value = "2.621.04"
$ws.Range("E14").Value = "  -0.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.69"
$ws.Range("E15").Value = "  +1.37%  "
$ws.Range("D16").Value = "2.257.79"
$ws.Range("E16").Value = "  -1.65%  "
$ws.Range("E17").Value = "  -1.08%  "
$ws.Range("D18").Value = "42.119.93"
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.26"
$ws.Range("E19").Value = "  -3.37%  "
$ws.Range("E20").Value = "  -1.45%  "
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.81"
$ws.Range("E22").Value = "  -0.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.21"
$ws.Range("E23").Value = "  -2.45%  "
$ws.Range("E24").Value = "  -0.81%  "
$ws.Range("E25").Value = "  +1.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.53"
$ws.Range("E27").Value = "  -2.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.78"
$ws.Range("E28").Value = "  +4.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.58"
$ws.Range("E29").Value = "  -0.44%  "
$ws.Range("E30").Value = "  +0.69%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "162.57"
$ws.Range("E31").Value = "  +0.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.25"
$ws.Range("E32").Value = "  -1.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.15"
$ws.Range("E34").Value = "  +2.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.61"
$ws.Range("E35").Value = "  +2.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0736"
$ws.Range("E36").Value = "  -2.49%  "
$ws.Range("E37").Value = "  -0.54%  "
$ws.Range("E38").Value = "  -3.67%  "
$ws.Range("E39").Value = "  -1.38%  "
$ws.Range("E40").Value = "  -1.10%  "
$ws.Range("E41").Value = "  -2.32%  "
$ws.Range("E42").Value = "  +3.00%  "
$ws.Range("D43").Value = "1.950.96"
$ws.Range("E43").Value = "  -3.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.00"
$ws.Range("E44").Value = "  -2.76%  "
$ws.Range("E45").Value = "  -1.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.01"
$ws.Range("E47").Value = "  -2.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.00"
$ws.Range("E48").Value = "  +1.03%  "
$ws.Range("D49").Value = "2.491.80"
$ws.Range("E49").Value = "  -0.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "92.14"
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.65"
$ws.Range("E51").Value = "  -1.93%  "
